$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Target: the paragraph holding the "_GoBack" bookmark (the very last
# paragraph in the body). We need to add a run of text before the bookmark
# and a run containing ")" after it, both formatted the same way as the rest
# of the body text: <w:rFonts w:hint="eastAsia"/>
# <w:lang w:val="en-US" w:eastAsia="zh-CN"/>  (this is already the paragraph
# mark's rPr on that paragraph, and also the rPr on the "}" run two
# paragraphs above it).
# ---------------------------------------------------------------------------

$bigText = "在多线程访问容器(List/Map)的情况下，为了避免同步加锁机制带了的串行化执行(牺牲效率),jdk新出了 CopyOnWriteList(先拷贝一份复制list,进行写操作，有并发读的时候访问原有容器，然后将指针指向这个新创建的list), ConcurrentHashMap(则是将map最多分16段，每段进行加锁，不同段可以异步操作"
$closeText = ")"

# Locate the bookmark and the paragraph that holds it (works no matter where
# in the body the bookmark lives).
$bookmark = $d.Bookmarks.Item("_GoBack")
$bmParaIndex = $bookmark.Range.Paragraphs.Item(1).Index

# A nearby run ("}" two paragraphs above the bookmark paragraph) that already
# carries the desired formatting; we borrow its FormattedText to stamp the
# same run properties onto our new runs further down.
$fmtParaRange = $d.Paragraphs.Item($bmParaIndex - 2).Range
$fmtSource = $d.Range($fmtParaRange.Start, $fmtParaRange.End - 1)

# Step 1: insert a placeholder run right AFTER the bookmark (this insertion
# order is what keeps the new text on the correct side of the bookmark).
$afterRange = $bookmark.Range
$afterRange.InsertAfter("TMP_AFTER")

# Step 2: insert a placeholder run right BEFORE the paragraph (and hence
# before the bookmark too).
$targetParagraph = $d.Paragraphs.Item($bmParaIndex)
$paraRange = $targetParagraph.Range
$beforeRange = $d.Range($paraRange.Start, $paraRange.Start)
$beforeRange.InsertBefore("TMP_BEFORE")

# Step 3: stamp formatting onto the "before" placeholder (collapses it down
# to the length of the source text), then overwrite its text with the real
# content -- this keeps the correct run properties while changing the text.
$paraRange2 = $d.Paragraphs.Item($bmParaIndex).Range
$beforeLen = "TMP_BEFORE".Length
$beforeFmtTarget = $d.Range($paraRange2.Start, $paraRange2.Start + $beforeLen)
$beforeFmtTarget.FormattedText = $fmtSource.FormattedText

$paraRange2b = $d.Paragraphs.Item($bmParaIndex).Range
$beforeTextTarget = $d.Range($paraRange2b.Start, $paraRange2b.Start + $fmtSource.Text.Length)
$beforeTextTarget.Text = $bigText

# Step 4: same two-phase fixup for the "after" placeholder (now offset by
# the length of the big text we just inserted).
$paraRange3 = $d.Paragraphs.Item($bmParaIndex).Range
$bigLen = $bigText.Length
$afterLen = "TMP_AFTER".Length
$afterFmtTarget = $d.Range($paraRange3.Start + $bigLen, $paraRange3.Start + $bigLen + $afterLen)
$afterFmtTarget.FormattedText = $fmtSource.FormattedText

$paraRange3b = $d.Paragraphs.Item($bmParaIndex).Range
$afterTextTarget = $d.Range($paraRange3b.Start + $bigLen, $paraRange3b.Start + $bigLen + $fmtSource.Text.Length)
$afterTextTarget.Text = $closeText

Write-Output "final: $($d.Paragraphs.Item($bmParaIndex).Range.Text)"
